# feat: add 2022-Q1 data
#
# The existing "总计" (totals) sheet (sheetId 6) is renamed to "2022-Q1" and its
# content is replaced with the fund-holdings detail table for the new quarter.
# A brand-new sheet named "总计" is inserted in its former place (sheetId 7) and
# rebuilt with the updated totals table (old rows shifted down by one, with a
# new first row for 2022-Q1).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Sheet shuffle: keep the physical file that currently backs "总计" (it will
#    become "2022-Q1"), and add a fresh sheet named "总计" right after it so
#    the tab order / sheetId allocation matches the target workbook.
#    Rename the old sheet away from "总计" FIRST so the name is free for the
#    brand-new sheet to take.
# ---------------------------------------------------------------------------
$totalsOld = $wb.Worksheets.Item("总计")
$totalsOld.Name = "2022-Q1"
$totalsNew = $wb.Worksheets.Add($null, $totalsOld)
$totalsNew.Name = "总计"

# ---------------------------------------------------------------------------
# 2. Build the new totals sheet ("总计") content:
#    header row + new 2022-Q1 row + the previous rows shifted down (index +1).
# ---------------------------------------------------------------------------
$totalsNew.Cells.Item(1, 2).Value = "日期"
$totalsNew.Cells.Item(1, 3).Value = "持有数量(只)"
$totalsNew.Cells.Item(1, 4).Value = "持有市值(亿元)"

$totalsRows = @(
    @(0, "2022-Q1", 9,  4.22),
    @(1, "2021-Q4", 14, 6.71),
    @(2, "2021-Q3", 10, 7.48),
    @(3, "2021-Q2", 9,  6.51),
    @(4, "2021-Q1", 22, 7.17),
    @(5, "2020-Q4", 14, 3.78)
)

$r = 2
foreach ($row in $totalsRows) {
    $totalsNew.Cells.Item($r, 1).Value = $row[0]
    $totalsNew.Cells.Item($r, 2).Value = $row[1]
    $totalsNew.Cells.Item($r, 3).Value = $row[2]
    $totalsNew.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Copy the bold/border/centered header-row & index-column formatting (style
# index 2 in the original workbook) from the renamed sheet, which still has
# it applied on its own header/index cells at this point.
$totalsOld.Range("B1").Copy()
$totalsNew.Range("B1:D1").PasteSpecial(-4122)
$totalsOld.Range("A2").Copy()
$totalsNew.Range("A2:A7").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Replace the (already renamed) "2022-Q1" sheet's content with the
#    fund-holdings detail table.
# ---------------------------------------------------------------------------
$fundSheet = $totalsOld
$fundSheet.Cells.Clear()

$fundSheet.Cells.Item(1, 2).Value = "基金代码"
$fundSheet.Cells.Item(1, 3).Value = "基金名称"
$fundSheet.Cells.Item(1, 4).Value = "基金规模"
$fundSheet.Cells.Item(1, 5).Value = "股票总仓位"
$fundSheet.Cells.Item(1, 6).Value = "仓位占比"
$fundSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$fundSheet.Cells.Item(1, 8).Value = "仓位排名"

# Columns D, E, F, G hold text-formatted numeric-looking values in the target
# workbook (t="inlineStr"), so pre-format those columns as Text before writing
# so the COM layer doesn't silently coerce them to numbers.
$fundSheet.Range("D2:G10").NumberFormat = "@"

$fundRows = @(
    @(0, "510880", "华泰柏瑞上证红利ETF",             "181.00", "97.22", "2.27", "4.1087", 10),
    @(1, "080005", "长盛量化红利混合",                 "2.66",   "69.88", "3.20", "0.0851", 3),
    @(2, "006729", "万家中证500指数增强A",             "1.04",   "93.64", "1.23", "0.0128", 8),
    @(3, "006730", "万家中证500指数增强C",             "0.61",   "93.64", "1.23", "0.0075", 8),
    @(4, "013802", "财通资管中证钢铁指数A",             "0.11",   "90.83", "2.53", "0.0028", 10),
    @(5, "006201", "景顺长城量化先锋混合",             "0.09",   "46.30", "0.82", "0.0007", 10),
    @(6, "001273", "民生加银新动力灵活配置混合A",       "0.04",   "68.44", "1.60", "0.0006", 9),
    @(7, "001274", "民生加银新动力灵活配置混合D",       "0.04",   "68.44", "1.60", "0.0006", 9),
    @(8, "013803", "财通资管中证钢铁指数C",             "0.02",   "90.83", "2.53", "0.0005", 10)
)

$r = 2
foreach ($row in $fundRows) {
    $fundSheet.Cells.Item($r, 1).Value = $row[0]
    $fundSheet.Cells.Item($r, 2).NumberFormat = "@"
    $fundSheet.Cells.Item($r, 2).Value = $row[1]
    $fundSheet.Cells.Item($r, 3).Value = $row[2]
    $fundSheet.Cells.Item($r, 4).Value = $row[3]
    $fundSheet.Cells.Item($r, 5).Value = $row[4]
    $fundSheet.Cells.Item($r, 6).Value = $row[5]
    $fundSheet.Cells.Item($r, 7).Value = $row[6]
    $fundSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Re-apply the header / index-column style (it was cleared above along with
# the rest of the old totals content) by pulling it back from the rebuilt
# "总计" sheet, which still carries style index 2 on its header + A column.
$totalsNew.Range("B1").Copy()
$fundSheet.Range("B1:H1").PasteSpecial(-4122)
$totalsNew.Range("A2").Copy()
$fundSheet.Range("A2:A10").PasteSpecial(-4122)

Write-Host "done"
